# Apply the "updated spring start indices" edit.
# Fills in the start-index readings (column F on sheet 1; columns C:F on
# sheet 2) that were measured/recorded for the remaining trials, and
# updates the sheet selections / active tab to match where the author
# left the cursor after entering the data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 2L2LT4ST_37T
$ws2 = $wb.Worksheets.Item(2)   # 2L3LT4ST_37T
$ws3 = $wb.Worksheets.Item(3)   # 2L4LT4ST_37T
$ws4 = $wb.Worksheets.Item(4)   # 2L5LT4ST_37T

# --- Sheet 1 (2L2LT4ST_37T): fill column F, rows 2-17 -----------------
$ws1.Range("F2").Value  = 0.43
$ws1.Range("F3").Value  = 0.42
$ws1.Range("F4").Value  = 0.26
$ws1.Range("F5").Value  = 0
$ws1.Range("F6").Value  = 0.09
$ws1.Range("F7").Value  = 0.2
$ws1.Range("F8").Value  = 0.17
$ws1.Range("F9").Value  = 0.17
$ws1.Range("F10").Value = 0.13
$ws1.Range("F11").Value = 0.22
$ws1.Range("F12").Value = 0.18
$ws1.Range("F13").Value = 0.14
$ws1.Range("F14").Value = 0.08
$ws1.Range("F15").Value = 0.2
$ws1.Range("F16").Value = 0.26
$ws1.Range("F17").Value = 0.17

# --- Sheet 2 (2L3LT4ST_37T): fill columns C:F, rows 2-17 ---------------
$ws2.Range("C2").Value = 0.16
$ws2.Range("D2").Value = 0.28
$ws2.Range("E2").Value = 0.37
$ws2.Range("F2").Value = 0.3

$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 0.36
$ws2.Range("F3").Value = 0.09

$ws2.Range("C4").Value = 0.21
$ws2.Range("D4").Value = 0.22
$ws2.Range("E4").Value = 0.26
$ws2.Range("F4").Value = 0.15

$ws2.Range("C5").Value = 0.12
$ws2.Range("D5").Value = 0.32
$ws2.Range("E5").Value = 0.02
$ws2.Range("F5").Value = 0

$ws2.Range("C6").Value = 0.12
$ws2.Range("D6").Value = 0.08
$ws2.Range("E6").Value = 0.2
$ws2.Range("F6").Value = 0.22

$ws2.Range("C7").Value = 0
$ws2.Range("D7").Value = 0.09
$ws2.Range("E7").Value = 0.18
$ws2.Range("F7").Value = 0.18

$ws2.Range("C8").Value = 0.09
$ws2.Range("D8").Value = 0.25
$ws2.Range("E8").Value = 0.28
$ws2.Range("F8").Value = 0.21

$ws2.Range("C9").Value = 0.12
$ws2.Range("D9").Value = 0.17
$ws2.Range("E9").Value = 0
$ws2.Range("F9").Value = 0.08

$ws2.Range("C10").Value = 0.11
$ws2.Range("D10").Value = 0
$ws2.Range("E10").Value = 0.05
$ws2.Range("F10").Value = 0.18

$ws2.Range("C11").Value = 0.17
$ws2.Range("D11").Value = 0.25
$ws2.Range("E11").Value = 0.15
$ws2.Range("F11").Value = 0.18

$ws2.Range("C12").Value = 0.03
$ws2.Range("D12").Value = 0.17
$ws2.Range("E12").Value = 0.22
$ws2.Range("F12").Value = 0.34

$ws2.Range("C13").Value = 0.17
$ws2.Range("D13").Value = 0.05
$ws2.Range("E13").Value = 0.2
$ws2.Range("F13").Value = 0.15

$ws2.Range("C14").Value = 0.11
$ws2.Range("D14").Value = 0.14
$ws2.Range("E14").Value = 0.07
$ws2.Range("F14").Value = 0.36

$ws2.Range("C15").Value = 0.25
$ws2.Range("D15").Value = 0.26
$ws2.Range("E15").Value = 0.18
$ws2.Range("F15").Value = 0.2

$ws2.Range("C16").Value = 0.27
$ws2.Range("D16").Value = 0.29
$ws2.Range("E16").Value = 0.2
$ws2.Range("F16").Value = 0.11

$ws2.Range("C17").Value = 0.06
$ws2.Range("D17").Value = 0.54
$ws2.Range("E17").Value = 0.05
$ws2.Range("F17").Value = 0.04

# --- Selections on each sheet (where the author left the cursor) ------
$ws1.Range("L30").Select() | Out-Null
$ws2.Range("J18").Select() | Out-Null

# Sheet 3 becomes the active / tabSelected sheet, cursor on E20.
$ws3.Activate() | Out-Null
$ws3.Range("E20").Select() | Out-Null

# Sheet 4's selection is unchanged (still C2), nothing to do there.
